$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 92.6317602074611
$ws.Range("C2").Value = 92.14546699825243
$ws.Range("D2").Value = 90.00696453232483
$ws.Range("E2").Value = 92.52426056004643
$ws.Range("B3").Value = 98.49830909693215
$ws.Range("C3").Value = 97.46376620635135
$ws.Range("D3").Value = 98.39621761805462
$ws.Range("E3").Value = 98.11122075468573
$ws.Range("B4").Value = 99.2263278418861
$ws.Range("C4").Value = 99.16072329807571
$ws.Range("D4").Value = 99.24460985935852
$ws.Range("E4").Value = 99.26863685356672
$ws.Range("B5").Value = 98.71196478810971
$ws.Range("C5").Value = 98.72885904061765
$ws.Range("D5").Value = 98.71467993312679
$ws.Range("E5").Value = 98.69680091019751
$ws.Range("B6").Value = 98.3120559264534
$ws.Range("C6").Value = 98.20502330782047
$ws.Range("D6").Value = 98.23638128680922
$ws.Range("E6").Value = 98.18385245356703
$ws.Range("B7").Value = 97.71966316302193
$ws.Range("C7").Value = 97.73529062291912
$ws.Range("D7").Value = 97.77542735833339
$ws.Range("E7").Value = 97.711875317159
$ws.Range("B8").Value = 97.25859628670493
$ws.Range("C8").Value = 97.22156674004633
$ws.Range("D8").Value = 97.26587962437384
$ws.Range("E8").Value = 97.21392745665615
$ws.Range("B9").Value = 95.89971400377856
$ws.Range("C9").Value = 95.88410862621126
$ws.Range("D9").Value = 95.88266255201073
$ws.Range("E9").Value = 95.91582352285025
